# "Switches excel file correction"
# Turn the populated switch-list sheet into a blank fill-in-the-blanks
# template: keep the header row, replace the first data row with
# placeholder prompts, and blank out the remaining sample rows while
# keeping their formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the placeholder / prompt row instead of real sample data,
# and needs extra height so the wrapped text fits.
$ws.Range("A2").Value = "<Enter Hostname>"
$ws.Range("B2").Value = "<Enter Device IP>"
$ws.Rows.Item(2).RowHeight = 28.8

# Remove the rest of the old sample switches, leaving the rows blank but
# still formatted (style carries over automatically with ClearContents).
$ws.Range("A3:B7").ClearContents()

# Column B no longer needs to best-fit the old IP-address samples; widen
# it to comfortably fit the new placeholder text instead.
$ws.Columns.Item(2).ColumnWidth = 17.5

# Leave the selection where the author ended up working.
$ws.Range("D6").Select()
